$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values are stored as literal text in the source data
# (t="inlineStr"), not as numbers, so force each target cell to Text format
# before writing the new numeric-looking string -- otherwise Excel would
# auto-convert it to a real number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "251.84"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.76"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.031"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05972"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.427"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.570"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.326"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.7989"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1490"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07890"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03356"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03055"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09295"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.560"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001683"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04774"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0006092"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006209"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005693"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001067"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.691"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.214"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3308"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1256"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0006477"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04434"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006997"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003601"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1067"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009122"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.002461"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005878"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.7855"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.09759"
$ws.Range("E49").Value = "48BOLOBOLO"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
